$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '59.784.33'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +4.21%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.026.32'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +4.00%  '
$ws.Range("E4").Value = '  +0.08%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '565.74'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +3.48%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '140.40'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +8.55%  '
$ws.Range("E7").Value = '  -0.16%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.521'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +2.18%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '3.017.25'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +4.02%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.134'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +6.94%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.28'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +11.41%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.461'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +3.69%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000233'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +6.38%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '34.14'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +4.40%  '
$ws.Range("E15").Value = '  +1.84%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.525.84'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +3.96%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '7.27'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +6.77%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.022.20'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +3.89%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '59.765.23'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +4.12%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '437.98'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +5.27%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.70'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +4.75%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.726'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +6.84%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.16'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +3.36%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '13.28'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +2.17%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '81.11'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.98%  '
$ws.Range("E26").Value = '  +0.20%  '
$ws.Range("E28").Value = '  +0.27%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.56'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +4.09%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.87'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +6.26%  '
$ws.Range("B31").Value = 'NEARProtocol'
$ws.Range("C31").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.30'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +6.87%  '
$ws.Range("B32").Value = 'EthereumClassic'
$ws.Range("C32").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '26.09'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +4.00%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.101'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +5.91%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0₃0791'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +17.94%  '
$ws.Range("E35").Value = '  +8.39%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.97'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +6.41%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.13'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +3.46%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '49.21'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +3.16%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '8.69'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.24%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.82'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +11.38%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '404.46'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +8.91%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0356'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +3.73%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.794.19'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +5.10%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.107'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.97%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.256'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +8.27%  '
$ws.Range("E46").Value = '  +0.01%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '123.42'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +1.47%  '
$ws.Range("B48").Value = 'Stellar'
$ws.Range("C48").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.111'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +2.01%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.03'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +3.89%  '
$ws.Range("B50").Value = 'Arweave'
$ws.Range("C50").Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '33.96'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +23.09%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '23.73'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +2.97%  '
